# Avance validaciones 1 al 10
# Update active selection on the sheet and widen a few columns to fit the
# expanded validator descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the selection to H10 (was B12)
$ws.Range("H10").Select()

# Widen columns B, C, D to fit the longer descriptions / codes.
# (ColumnWidth is expressed in the workbook's "characters" unit and gets
# quantized to whole pixels on save, so we pick the nearest value that
# rounds back to the target stored width.)
$ws.Columns.Item(2).ColumnWidth = 75.66666666666667
$ws.Columns.Item(3).ColumnWidth = 9
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
